$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Move the "Type" header/description from row 4 to row 5 so it lines up
# with the "Type" label in column A (row 5).
$typeHeader = $ws.Range("F4").Value2
$typeDesc = $ws.Range("G4").Value2
$ws.Range("F4:G4").ClearContents()
$ws.Range("F5").Value2 = $typeHeader
$ws.Range("G5").Value2 = $typeDesc

# Remove the "Sectors" entry (row 10) from the controlled-vocab list.
$ws.Range("A10").ClearContents()

# Update the active selection to match the saved view state.
$ws.Range("F2").Select()
